$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "321.60"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.33%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "49.14"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "11.24%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.316"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.48%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08065"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.69%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.604"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.23%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.352"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "29.00%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.80%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1280"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.41%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1968"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.38%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09642"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.17%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04722"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "13.30%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.28%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001319"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "3.22%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.44%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005903"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.43%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.345"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.04%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.444"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.86%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3509"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "4.98%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.017"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.07%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1363"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.54%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3093"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.67%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001315"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004307"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.45%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001350"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.02%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003539"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02727"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "8.85%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05984"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "13.04%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01086"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "83.03%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.008027"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.80%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.63%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007897"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "7.43%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008663"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "14.74%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3500"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "16.17%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006898"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.34%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.04%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "37.49%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003999"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-4.80%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.04%"
